# Update the constant AO value (rows 3:18) used in the SMP_3PP equilibrium
# equations, then let Excel recalculate the dependent AR/AS/AT/AU formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AO3:AO18").Value = 317303.66712625924

$excel.CalculateFullRebuild()
